$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Habitats and specific countries"
$ws.Range("D1").Value = "Size (length or weight)"

$ws.Range("B3").Value = "Wetlands (inland) - Permanent Rivers/Streams/Creeks (includes waterfalls), Wetlands (inland) - Permanent Freshwater Marshes/Pools (under 8ha), Texas, Mexico, United States"

$ws.Range("B4").Value = "Marine Neritic - Seagrass (Submerged), Marine Neritic - Subtidal Loose Rock/pebble/gravel, Marine Neritic - Subtidal Rock and Rocky Reefs, Norway, Spain, Portugal, France, Germany, Italy, Greece, Türkiye, Syrian Arab Republic, Lebanon, Tunisia, Morocco, Malta, Jersey, Guernsey, Gibraltar, Denmark, Croatia, Bulgaria, Belgium, Albania, Monaco, Montenegro, Netherlands, Romania, Slovenia, Algeria, Egypt, Libya"
$ws.Range("C4").Value = "small gastropods, sea urchins, worms, shrimps, isopods, amphipods"

$ws.Range("B5").Value = "Savanna - Moist, Forest - Subtropical/Tropical Moist Lowland, Forest - Subtropical/Tropical Swamp, Brazil, Amazonas, Colombia, Ecuador, Peru, Venezuela"
$ws.Range("C5").Value = "fruits, young leaves, petioles, flowers"

$ws.Range("B6").Value = "Forest - Temperate, Wetlands (inland) - Bogs, Marshes, Swamps, Fens, Peatlands, Wetlands (inland) - Shrub Dominated Wetlands, Grassland - Temperate, Shrubland - Temperate, Artificial/Terrestrial - Urban Areas, Artificial/Terrestrial - Rural Gardens, Artificial/Terrestrial - Pastureland, Artificial/Terrestrial - Arable Land, Wisconsin, West Virginia, Virginia, Vermont, Texas, Tennessee, South Dakota, South Carolina, Rhode Island, Pennsylvania, Oklahoma, Ohio, North Dakota, North Carolina, New York, New Jersey, New Hampshire, Nebraska, Montana, Missouri, Mississippi, Minnesota, District of Columbia, Delaware, Connecticut,"
